$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (Arkansas)
$ws.Range("C10").Value = 35246
$ws.Range("D10").Value = 380
$ws.Range("E10").Value = 7521
$ws.Range("F10").Value = 100
$ws.Range("G10").Value = 24.57
$ws.Range("H10").Value = 26.53
$ws.Range("K10").Value = 30608
$ws.Range("L10").Value = 377

# Row 41 (Iowa)
$ws.Range("D41").Value = 812
$ws.Range("E41").Value = 3287
$ws.Range("H41").Value = 4.68
